$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-03 10:56:08"

$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

for ($row = 2; $row -le 26; $row++) {
    $addr = "AA" + $row
    $ws2.Range($addr).Value = $newTimestamp
    $ws3.Range($addr).Value = $newTimestamp
}
